$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 8 rows (19-26) that held the frequency-related block values:
# GnssPoseSimulink, PointsRawFloat32, ImageRaw, ClockFrequency, SimulinkState,
# CurrentVelocity, PoseOtherCar and CurrentPose. Deleting the entire rows shifts
# the remaining rows (percent_reflecting_sfc, R) up so they become rows 19-20.
$ws.Range("A19:B26").EntireRow.Delete() | Out-Null

# Reflect the view/selection state left behind after performing the deletion.
$excel.ActiveWindow.ScrollRow = 7 | Out-Null
$ws.Range("A19:XFD26").Select() | Out-Null
